$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 58

$ws.Cells.Item($row, 1).Value = "2024-09-27 00:00:00"
$ws.Cells.Item($row, 2).Value = 75550
$ws.Cells.Item($row, 3).Value = 10762.57
$ws.Cells.Item($row, 4).Value = 9524.4
$ws.Cells.Item($row, 5).Value = 7.0113
